$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header cell in H1, matching the formatting of the existing
# header row (copy G1's format onto H1, then set the text).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data values for the new "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
